$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.707.21'

$ws.Range('D3').Value = '1.889.69'
$ws.Range('E3').Value = '  +1.39%  '

$ws.Range('E4').Value = '  -0.85%  '

$ws.Range('D5').Value = "'313.55"
$ws.Range('E5').Value = '  +0.11%  '

$ws.Range('E6').Value = '  -0.90%  '

$ws.Range('D7').Value = "'0.4803"
$ws.Range('E7').Value = '  +0.43%  '

$ws.Range('D8').Value = "'0.3779"
$ws.Range('E8').Value = '  -0.67%  '

$ws.Range('D9').Value = "'0.07319"
$ws.Range('E9').Value = '  -0.26%  '

$ws.Range('D10').Value = "'0.9182"
$ws.Range('E10').Value = '  -1.57%  '

$ws.Range('D11').Value = "'20.38"
$ws.Range('E11').Value = '  -1.82%  '

$ws.Range('D12').Value = "'0.07696"
$ws.Range('E12').Value = '  -1.10%  '

$ws.Range('D13').Value = '1.901.68'
$ws.Range('E13').Value = '  +1.89%  '

$ws.Range('D14').Value = "'5.463"
$ws.Range('E14').Value = '  +0.40%  '

$ws.Range('D15').Value = "'6.576"
$ws.Range('E15').Value = '  +0.09%  '

$ws.Range('D16').Value = "'90.79"
$ws.Range('E16').Value = '  +0.55%  '

$ws.Range('D17').Value = "'1.004"
$ws.Range('E17').Value = '  -0.86%  '

$ws.Range('D18').Value = "'0.000008799"

$ws.Range('E19').Value = '  -0.80%  '

$ws.Range('D20').Value = '27.755.89'
$ws.Range('E20').Value = '  +0.16%  '

$ws.Range('D21').Value = "'14.47"
$ws.Range('E21').Value = '  -1.57%  '

$ws.Range('D22').Value = "'5.113"
$ws.Range('E22').Value = '  +0.26%  '

$ws.Range('D23').Value = '2.149.83'
$ws.Range('E23').Value = '  +0.81%  '

$ws.Range('D24').Value = "'10.80"
$ws.Range('E24').Value = '  +0.76%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = "'1.900"
$ws.Range('E25').Value = '  -1.87%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = "'153.75"
$ws.Range('E26').Value = '  -1.42%  '

$ws.Range('D27').Value = "'18.35"
$ws.Range('E27').Value = '  -0.80%  '

$ws.Range('D28').Value = "'2.107"
$ws.Range('E28').Value = '  +4.22%  '

$ws.Range('D29').Value = "'116.12"
$ws.Range('E29').Value = '  +0.61%  '

$ws.Range('D30').Value = "'4.911"
$ws.Range('E30').Value = '  -1.09%  '

$ws.Range('D31').Value = "'0.08917"
$ws.Range('E31').Value = '  +0.47%  '

$ws.Range('D32').Value = "'3.155"
$ws.Range('E32').Value = '  -5.26%  '

$ws.Range('D33').Value = "'1.229"
$ws.Range('E33').Value = '  +1.94%  '

$ws.Range('D34').Value = "'0.7598"
$ws.Range('E34').Value = '  +0.71%  '

$ws.Range('D35').Value = "'4.623"
$ws.Range('E35').Value = '  +0.70%  '

$ws.Range('D36').Value = "'0.02032"
$ws.Range('E36').Value = '  -0.31%  '

$ws.Range('D37').Value = "'2.520"
$ws.Range('E37').Value = '  -6.68%  '

$ws.Range('D38').Value = "'1.090"
$ws.Range('E38').Value = '  -3.02%  '

$ws.Range('D39').Value = "'0.05249"
$ws.Range('E39').Value = '  -1.91%  '

$ws.Range('D40').Value = "'2.974"
$ws.Range('E40').Value = '  -0.24%  '

$ws.Range('D41').Value = "'0.5422"
$ws.Range('E41').Value = '  -4.60%  '

$ws.Range('D42').Value = "'6.954"
$ws.Range('E42').Value = '  -1.11%  '

$ws.Range('D43').Value = "'0.1513"
$ws.Range('E43').Value = '  -0.86%  '

$ws.Range('D44').Value = "'8.296"
$ws.Range('E44').Value = '  -2.60%  '

$ws.Range('D45').Value = "'109.36"
$ws.Range('E45').Value = '  +4.14%  '

$ws.Range('D46').Value = "'10.65"
$ws.Range('E46').Value = '  -0.61%  '

$ws.Range('D47').Value = "'0.4770"
$ws.Range('E47').Value = '  -2.30%  '

$ws.Range('D48').Value = "'1.002"
$ws.Range('E48').Value = '  -0.98%  '

$ws.Range('D49').Value = "'1.634"
$ws.Range('E49').Value = '  -1.73%  '

$ws.Range('D50').Value = "'67.48"
$ws.Range('E50').Value = '  -0.08%  '

$ws.Range('D51').Value = "'0.06058"
$ws.Range('E51').Value = '  -0.68%  '

